$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain value assignments (unambiguous text; Excel will keep these as text)
$ws.Range("D2").Value = "57.842.98"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.455.95"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D9").Value = "2.464.02"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "2.892.54"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "57.790.73"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "2.459.74"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -4.83%  "
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "1.720.07"
$ws.Range("E51").Value = "  -1.78%  "

# Ambiguous numeric-looking text values: force Text format first so Excel
# stores them as strings (matching the original inline-string cell type)
# instead of auto-converting to numbers.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "516.94"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "132.30"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0972"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.27"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.86"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "317.72"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.68"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "64.40"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.27"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "168.33"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "17.94"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "36.29"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "270.26"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.95"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.587"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "123.40"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0903"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0483"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "16.66"
